$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13-23 shift down to 14-24.
$ws.Rows(13).Insert()

# Populate the newly inserted row 13 with the new weekly price record.
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 44790
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino dulce"
$ws.Range("H13").Value = "Cultivar IV Región"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17500
$ws.Range("N13").Value = "$/bandeja 18 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 972
$ws.Range("Q13").Value = 18
$ws.Range("R13").Value = "Hortaliza"
